{"js": "// Fix the placeholder fields in the \"Tempat & Tanggal Lahir\" line so that\n// they refer to the maker's (pembuat) birth place / date instead of the\n// generic ones, matching the {xxxPembuat} naming convention used by the\n// rest of the template (e.g. {namaPembuat}, {alamatPembuat}).\n//\n//   {tempatLahir}  ->  {tempatLahirPembuat}\n//   {ttl}          ->  {ttlPembuat}\n\nconst body = context.document.body;\n\n// --- Fix {tempatLahir} -> {tempatLahirPembuat} -----------------------------\nlet results = body.search(\"{tempatLahir}\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"{tempatLahirPembuat}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Fix {ttl} -> {ttlPembuat} ---------------------------------------------\nresults = body.search(\"{ttl}\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"{ttlPembuat}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix the placeholder fields in the \"Tempat & Tanggal Lahir\" line so that\n# they refer to the maker's (pembuat) birth place / date instead of the\n# generic ones, matching the {xxxPembuat} naming convention used by the\n# rest of the template (e.g. {namaPembuat}, {alamatPembuat}).\n#\n#   {tempatLahir}  ->  {tempatLahirPembuat}\n#   {ttl}          ->  {ttlPembuat}\n\n$d = $word.ActiveDocument\n\n# --- Fix {tempatLahir} -> {tempatLahirPembuat} ------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"{tempatLahir}\"\n$find.Replacement.Text = \"{tempatLahirPembuat}\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) # wdReplaceAll\n\n# --- Fix {ttl} -> {ttlPembuat} ----------------------------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"{ttl}\"\n$find2.Replacement.Text = \"{ttlPembuat}\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
